$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 2.6639427081334102
$ws.Range("B2").Value = 4.3728187783424124
$ws.Range("C2").Value = 211319

$ws.Range("A3").Value = 2.7070197476837885
$ws.Range("B3").Value = 3.0009077762335861
$ws.Range("C3").Value = 139975

$ws.Range("A6").Value = 2.699344872651511
$ws.Range("B6").Value = 3.0957296363440197
$ws.Range("C6").Value = 192572
